$d = $word.ActiveDocument

function Set-ParaXML($para, $innerXml) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
           $innerXml + `
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

$dash = [char]0x2013

# --- 1. Merge split runs / remove proofErr wrappers for paragraphs whose text is unchanged
#        but whose run/proofErr structure collapses into a single run. Paragraph indices
#        below refer to the ORIGINAL (pre-edit) document and remain stable because each
#        replacement keeps the paragraph count the same (1 paragraph in, 1 paragraph out).

Set-ParaXML $d.Paragraphs.Item(1) ('<w:p><w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Assets from Prototype2 as it is a 2D game - will change drawings later</w:t></w:r></w:p>')

Set-ParaXML $d.Paragraphs.Item(5) '<w:p><w:r><w:t>Instructions screen</w:t></w:r></w:p>'

Set-ParaXML $d.Paragraphs.Item(8) '<w:p><w:r><w:t>Create game characters</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(9) '<w:p><w:r><w:t>Create background</w:t></w:r></w:p>'

Set-ParaXML $d.Paragraphs.Item(11) '<w:p><w:r><w:t>Create rocks</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(12) '<w:p><w:r><w:t>Create buckets</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(13) '<w:p><w:r><w:t>Create stone</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(14) '<w:p><w:r><w:t>Create rope</w:t></w:r></w:p>'

Set-ParaXML $d.Paragraphs.Item(19) '<w:p><w:r><w:t>Pause icon</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(20) '<w:p><w:r><w:t>Pause Screen - resume, restart, exit</w:t></w:r></w:p>'

Set-ParaXML $d.Paragraphs.Item(23) '<w:p><w:r><w:t>Create running plane and background plane</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(24) '<w:p><w:r><w:t>Create player</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(25) '<w:p><w:r><w:t>Create player controller - how to move</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(26) '<w:p><w:r><w:t>Create spawner - obstacles and well pieces</w:t></w:r></w:p>'

Set-ParaXML $d.Paragraphs.Item(28) '<w:p><w:r><w:t>Respawn/ change level</w:t></w:r></w:p>'
Set-ParaXML $d.Paragraphs.Item(29) '<w:p><w:r><w:t>Increase difficulty for each level</w:t></w:r></w:p>'

# --- 2. Insert the five new "Main menu" to-do items between "Start Screen" (item 4) and
#        "Instructions screen" (item 5).

$newParasXml = '<w:p><w:r><w:t xml:space="preserve">Options Button </w:t></w:r><w:r><w:t>' + $dash + '</w:t></w:r><w:r><w:t xml:space="preserve"> sound</w:t></w:r></w:p>' + `
               '<w:p><w:r><w:t>Sound effects</w:t></w:r></w:p>' + `
               '<w:p><w:r><w:t>Music</w:t></w:r></w:p>' + `
               '<w:p><w:r><w:t>Quit Button</w:t></w:r></w:p>' + `
               '<w:p><w:r><w:t>Play button</w:t></w:r></w:p>'

$startScreen = $d.Paragraphs.Item(4)
$startScreen.Range.InsertParagraphAfter()
$newEmptyPara = $d.Paragraphs.Item(5)
Set-ParaXML $newEmptyPara $newParasXml
